# Sprint 3 Burndown Chart - apply "Updated burndown chart and removed remnants of the past"
#
# Summary of the change being reproduced:
#   1. A new data row (row 11) is appended under the existing burndown table
#      on Sheet1.
#   2. The embedded bar chart's category/value ranges grow from $2:$10 to
#      $2:$11 to pick the new row up.
#   3. The chart title / axis titles get small wording tweaks.
#   4. The chart's on-sheet anchor (size/position) is nudged.
#   5. The saved cursor position moves to D22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the new burndown row. Copy the previous (last) row first so the new
#    cells inherit the same number formats/styles (date format on A, percent
#    format on C) instead of minting new style entries, then overwrite with
#    the real values.
# ---------------------------------------------------------------------------
$ws.Range("A10:C10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 44631
$ws.Range("B11").Value = "Finished basement as well as miscellainous items"
$ws.Range("C11").Value = 0

# ---------------------------------------------------------------------------
# 2. Grow the chart series ranges to include the new row.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$A`$11,Sheet1!`$C`$2:`$C`$11,1)"

# ---------------------------------------------------------------------------
# 3. Wording tweaks on the chart/axis titles.
# ---------------------------------------------------------------------------
$chart.ChartTitle.Text = "Sprint 3 Burndown Chart"
$chart.Axes(1).AxisTitle.Text = "Date"
$chart.Axes(2).AxisTitle.Text = "Percent of tasks done"

# ---------------------------------------------------------------------------
# 4. Resize / reposition the chart's anchor on the sheet (same size, nudged
#    up and to the right slightly).
# ---------------------------------------------------------------------------
$co.Left = 659.9111328125
$co.Top = 2.62496062992126

# ---------------------------------------------------------------------------
# 5. Leave the cursor where the author last left it.
# ---------------------------------------------------------------------------
$ws.Range("D22").Select()
